$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1325
$ws1.Range("F3").Value = 1209
$ws1.Range("F4").Value = 14572
$ws1.Range("F5").Value = 17562
$ws1.Range("F7").Value = 67
$ws1.Range("F9").Value = 212
$ws1.Range("F14").Value = 6
$ws1.Range("F15").Value = 42
$ws1.Range("F16").Value = 34
$ws1.Range("F17").Value = 139
$ws1.Range("F19").Value = 1324
$ws1.Range("F20").Value = 147
$ws1.Range("F22").Value = 63
$ws1.Range("F23").Value = 210
$ws1.Range("F24").Value = 7214
$ws1.Range("F25").Value = 978
$ws1.Range("F26").Value = 1
$ws1.Range("F27").Value = 35
$ws1.Range("F28").Value = 1168
$ws1.Range("F29").Value = 10
$ws1.Range("F30").Value = 5855
$ws1.Range("F31").Value = 62
$ws1.Range("F32").Value = 43
$ws1.Range("F33").Value = 133
$ws1.Range("F35").Value = 222
$ws1.Range("F36").Value = 5075
$ws1.Range("F37").Value = 20

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1325
$ws4.Range("F3").Value = 1209
$ws4.Range("F4").Value = 14572
$ws4.Range("F5").Value = 17562
$ws4.Range("F7").Value = 67
$ws4.Range("F9").Value = 212
$ws4.Range("F14").Value = 6
$ws4.Range("F15").Value = 42
$ws4.Range("F16").Value = 34
$ws4.Range("F17").Value = 139
$ws4.Range("F19").Value = 1324
$ws4.Range("F20").Value = 147
$ws4.Range("F23").Value = 63
$ws4.Range("F24").Value = 210
$ws4.Range("F25").Value = 7214
$ws4.Range("F27").Value = 1
$ws4.Range("F28").Value = 35
$ws4.Range("F29").Value = 1168
$ws4.Range("F32").Value = 5855
$ws4.Range("F33").Value = 62
$ws4.Range("F34").Value = 43
$ws4.Range("F35").Value = 133
$ws4.Range("F37").Value = 222
$ws4.Range("F38").Value = 5075
$ws4.Range("F39").Value = 20
